# Weekly update: insert a new week's Berenjena (eggplant) record.
# This shifts the existing rows 66-100 down to 67-101 (dimension grows
# from A1:R100 to A1:R101) and fills the newly opened row 66 with the
# new week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 66, pushing rows 66-100
# down to 67-101.
$ws.Rows.Item(66).Insert()

# Fill the newly inserted row 66 with the new weekly record.
$ws.Range("A66").Value = 5
$ws.Range("B66").Value = "Macroferia Regional de Talca"
$ws.Range("C66").Value = "Maule"
$ws.Range("D66").Value = 44603
$ws.Range("E66").Value = 7
$ws.Range("F66").Value = 100112001
$ws.Range("G66").Value = "Berenjena"
$ws.Range("H66").Value = "Sin especificar"
$ws.Range("I66").Value = "Primera"
$ws.Range("J66").Value = 200
$ws.Range("K66").Value = 7000
$ws.Range("L66").Value = 7000
$ws.Range("M66").Value = 7000
$ws.Range("N66").Value = "$/caja 50 unidades"
$ws.Range("O66").Value = "Región del Maule"
$ws.Range("P66").Value = 140
$ws.Range("Q66").Value = 50
$ws.Range("R66").Value = "Hortaliza"
